$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the email body (F2) to reference where the job applications came from.
# Setting this first means the new shared string for the body text is appended
# to the shared-string table before the new recipient-email string.
$ws.Range("F2").Value = "Job application(s) submitted from :`n{0}`nThis email is sent by an UiPath Orchestrator bot.`nRegards"

# The recipient email (B2) changes to a new address. Remove the old hyperlink
# that pointed at the previous (hotmail) address, then write the new value.
# A bulk Hyperlinks.Delete() on the sheet clears every hyperlink, so capture
# the two that must survive and re-create them afterwards.
$eHyperlink = $ws.Range("E2").Hyperlinks.Item(1)
$eAddress = $eHyperlink.Address
$bHyperlink = $ws.Range("B3:B6").Hyperlinks.Item(1)
$bAddress = $bHyperlink.Address
$bDisplay = $bHyperlink.TextToDisplay

$ws.Hyperlinks.Delete()

$ws.Range("B2").Value = " ibnur@raceacademy.com.sg "

$ws.Hyperlinks.Add($ws.Range("E2"), $eAddress) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3:B6"), $bAddress, "", "", $bDisplay) | Out-Null

$ws.Range("C6").Select()
